$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1757322175732217
$ws.Range("C2").Value = 0.5648535564853556
$ws.Range("J2").Value = 0.008368200836820083
$ws.Range("P2").Value = 0.1297071129707113
$ws.Range("S2").Value = 0.1213389121338912

# Row 3
$ws.Range("B3").Value = 0.007194244604316547
$ws.Range("C3").Value = 0.03597122302158273
$ws.Range("J3").Value = 0.02877697841726619
$ws.Range("P3").Value = 0.7482014388489209
$ws.Range("S3").Value = 0.1798561151079137

# Row 4
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.7586206896551724
$ws.Range("S4").Value = 0.1724137931034483

# Row 6
$ws.Range("B6").Value = 0.05194805194805195
$ws.Range("D6").Value = 0.008658008658008658
$ws.Range("F6").Value = 0.05194805194805195
$ws.Range("J6").Value = 0.2943722943722944
$ws.Range("O6").Value = 0.01298701298701299
$ws.Range("Q6").Value = 0.1904761904761905
$ws.Range("R6").Value = 0.06493506493506493
$ws.Range("S6").Value = 0.3246753246753247

# Row 7
$ws.Range("B7").Value = 0.1058201058201058
$ws.Range("D7").Value = 0.02116402116402116
$ws.Range("F7").Value = 0.06349206349206349
$ws.Range("J7").Value = 0.1534391534391534
$ws.Range("O7").Value = 0.02645502645502645
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("R7").Value = 0.08994708994708994
$ws.Range("S7").Value = 0.3968253968253968

# Row 8
$ws.Range("B8").Value = 0.08836206896551724
$ws.Range("D8").Value = 0.01724137931034483
$ws.Range("F8").Value = 0.09482758620689655
$ws.Range("J8").Value = 0.08620689655172414
$ws.Range("O8").Value = 0.01939655172413793
$ws.Range("Q8").Value = 0.2025862068965517
$ws.Range("R8").Value = 0.06896551724137931
$ws.Range("S8").Value = 0.4224137931034483

# Row 9
$ws.Range("B9").Value = 0.08866995073891626
$ws.Range("D9").Value = 0.009852216748768473
$ws.Range("F9").Value = 0.04926108374384237
$ws.Range("J9").Value = 0.1379310344827586
$ws.Range("O9").Value = 0.02955665024630542
$ws.Range("Q9").Value = 0.2167487684729064
$ws.Range("R9").Value = 0.0541871921182266
$ws.Range("S9").Value = 0.4137931034482759

# Row 10
$ws.Range("B10").Value = 0.09803921568627451
$ws.Range("D10").Value = 0.0130718954248366
$ws.Range("E10").Value = 0.002801120448179272
$ws.Range("F10").Value = 0.08123249299719888
$ws.Range("J10").Value = 0.1092436974789916
$ws.Range("O10").Value = 0.0130718954248366
$ws.Range("Q10").Value = 0.2091503267973856
$ws.Range("R10").Value = 0.0718954248366013
$ws.Range("S10").Value = 0.4014939309056956

# Row 11
$ws.Range("G11").Value = 0.146551724137931
$ws.Range("J11").Value = 0.06465517241379311
$ws.Range("K11").Value = 0.1508620689655172
$ws.Range("L11").Value = 0.6293103448275862
$ws.Range("S11").Value = 0.008620689655172414

# Row 12
$ws.Range("G12").Value = 0.757396449704142
$ws.Range("J12").Value = 0.1893491124260355
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.01183431952662722
$ws.Range("S12").Value = 0.03550295857988166

# Row 13
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2653061224489796
$ws.Range("S13").Value = 0.02040816326530612

# Row 15
$ws.Range("F15").Value = 0.02463054187192118
$ws.Range("H15").Value = 0.187192118226601
$ws.Range("I15").Value = 0.07881773399014778
$ws.Range("J15").Value = 0.3399014778325123
$ws.Range("K15").Value = 0.06896551724137931
$ws.Range("M15").Value = 0.02463054187192118
$ws.Range("O15").Value = 0.03448275862068965
$ws.Range("S15").Value = 0.2413793103448276

# Row 16
$ws.Range("F16").Value = 0.01298701298701299
$ws.Range("H16").Value = 0.2337662337662338
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3506493506493507
$ws.Range("K16").Value = 0.07792207792207792
$ws.Range("M16").Value = 0.05844155844155844
$ws.Range("O16").Value = 0.09740259740259741
$ws.Range("S16").Value = 0.07792207792207792

# Row 17
$ws.Range("F17").Value = 0.01851851851851852
$ws.Range("H17").Value = 0.1990740740740741
$ws.Range("I17").Value = 0.1041666666666667
$ws.Range("J17").Value = 0.3726851851851852
$ws.Range("K17").Value = 0.1041666666666667
$ws.Range("M17").Value = 0.02314814814814815
$ws.Range("O17").Value = 0.06944444444444445
$ws.Range("S17").Value = 0.1087962962962963

# Row 18
$ws.Range("F18").Value = 0.01333333333333333
$ws.Range("H18").Value = 0.2466666666666667
$ws.Range("I18").Value = 0.04
$ws.Range("J18").Value = 0.36
$ws.Range("K18").Value = 0.1666666666666667
$ws.Range("M18").Value = 0.01333333333333333
$ws.Range("N18").Value = 0.006666666666666667
$ws.Range("O18").Value = 0.08666666666666667
$ws.Range("S18").Value = 0.06666666666666667

# Row 19
$ws.Range("F19").Value = 0.02300884955752212
$ws.Range("H19").Value = 0.2292035398230088
$ws.Range("I19").Value = 0.1053097345132743
$ws.Range("J19").Value = 0.3699115044247788
$ws.Range("K19").Value = 0.08141592920353982
$ws.Range("M19").Value = 0.01946902654867257
$ws.Range("O19").Value = 0.07256637168141593
$ws.Range("S19").Value = 0.09911504424778761
